$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("B2").Value = 0.21671826625387
$ws.Range("C2").Value = 0.5015479876160991
$ws.Range("J2").Value = 0.01857585139318885
$ws.Range("P2").Value = 0.1609907120743034
$ws.Range("S2").Value = 0.1021671826625387

# Row 3
$ws.Range("B3").Value = 0.02325581395348837
$ws.Range("C3").Value = 0.05232558139534884
$ws.Range("J3").Value = 0.04651162790697674
$ws.Range("P3").Value = 0.7034883720930233
$ws.Range("S3").Value = 0.1744186046511628

# Row 4
$ws.Range("J4").Value = 0.05263157894736842
$ws.Range("P4").Value = 0.6842105263157895
$ws.Range("S4").Value = 0.2631578947368421

# Row 6
$ws.Range("B6").Value = 0.03964757709251102
$ws.Range("D6").Value = 0.00881057268722467
$ws.Range("F6").Value = 0.04845814977973568
$ws.Range("J6").Value = 0.2378854625550661
$ws.Range("O6").Value = 0.013215859030837
$ws.Range("Q6").Value = 0.1806167400881057
$ws.Range("R6").Value = 0.06167400881057269
$ws.Range("S6").Value = 0.4096916299559472

# Row 7
$ws.Range("B7").Value = 0.1064814814814815
$ws.Range("D7").Value = 0.01388888888888889
$ws.Range("F7").Value = 0.04166666666666666
$ws.Range("J7").Value = 0.1851851851851852
$ws.Range("O7").Value = 0.009259259259259259
$ws.Range("Q7").Value = 0.1851851851851852
$ws.Range("R7").Value = 0.05555555555555555
$ws.Range("S7").Value = 0.4027777777777778

# Row 8
$ws.Range("B8").Value = 0.07592190889370933
$ws.Range("D8").Value = 0.006507592190889371
$ws.Range("F8").Value = 0.06724511930585683
$ws.Range("J8").Value = 0.1366594360086768
$ws.Range("O8").Value = 0.008676789587852495
$ws.Range("Q8").Value = 0.1626898047722343
$ws.Range("R8").Value = 0.08459869848156182
$ws.Range("S8").Value = 0.4577006507592191

# Row 9
$ws.Range("B9").Value = 0.04784688995215311
$ws.Range("D9").Value = 0.02392344497607655
$ws.Range("F9").Value = 0.03349282296650718
$ws.Range("J9").Value = 0.1052631578947368
$ws.Range("O9").Value = 0.01435406698564593
$ws.Range("Q9").Value = 0.1531100478468899
$ws.Range("R9").Value = 0.04784688995215311
$ws.Range("S9").Value = 0.5741626794258373

# Row 10
$ws.Range("B10").Value = 0.1115564462257849
$ws.Range("D10").Value = 0.01736806947227789
$ws.Range("E10").Value = 0.002672010688042752
$ws.Range("F10").Value = 0.072812291249165
$ws.Range("J10").Value = 0.1362725450901804
$ws.Range("O10").Value = 0.01269205076820307
$ws.Range("Q10").Value = 0.2064128256513026
$ws.Range("R10").Value = 0.0614562458249833
$ws.Range("S10").Value = 0.3787575150300601

# Row 11
$ws.Range("F11").Value = 0.002976190476190476
$ws.Range("G11").Value = 0.1577380952380952
$ws.Range("J11").Value = 0.09226190476190477
$ws.Range("K11").Value = 0.1875
$ws.Range("L11").Value = 0.5535714285714286
$ws.Range("S11").Value = 0.005952380952380952

# Row 12
$ws.Range("G12").Value = 0.7164948453608248
$ws.Range("J12").Value = 0.2216494845360825
$ws.Range("L12").Value = 0.02061855670103093
$ws.Range("S12").Value = 0.04123711340206185

# Row 13
$ws.Range("G13").Value = 0.6923076923076923
$ws.Range("J13").Value = 0.3076923076923077

# Row 15
$ws.Range("F15").Value = 0.0170940170940171
$ws.Range("H15").Value = 0.1495726495726496
$ws.Range("I15").Value = 0.07264957264957266
$ws.Range("J15").Value = 0.4017094017094017
$ws.Range("K15").Value = 0.07692307692307693
$ws.Range("M15").Value = 0.008547008547008548
$ws.Range("O15").Value = 0.04700854700854701
$ws.Range("S15").Value = 0.2264957264957265

# Row 16
$ws.Range("F16").Value = 0.015625
$ws.Range("H16").Value = 0.140625
$ws.Range("I16").Value = 0.09375
$ws.Range("J16").Value = 0.453125
$ws.Range("K16").Value = 0.1145833333333333
$ws.Range("M16").Value = 0.03125
$ws.Range("O16").Value = 0.02604166666666667
$ws.Range("S16").Value = 0.125

# Row 17
$ws.Range("F17").Value = 0.02244897959183673
$ws.Range("H17").Value = 0.1816326530612245
$ws.Range("I17").Value = 0.08163265306122448
$ws.Range("J17").Value = 0.4877551020408163
$ws.Range("K17").Value = 0.07551020408163266
$ws.Range("M17").Value = 0.00816326530612245
$ws.Range("O17").Value = 0.05510204081632653
$ws.Range("S17").Value = 0.08775510204081632

# Row 18
$ws.Range("F18").Value = 0.01796407185628742
$ws.Range("H18").Value = 0.1676646706586826
$ws.Range("I18").Value = 0.1137724550898204
$ws.Range("J18").Value = 0.4011976047904192
$ws.Range("K18").Value = 0.1017964071856287
$ws.Range("M18").Value = 0.01796407185628742
$ws.Range("N18").Value = 0.005988023952095809
$ws.Range("O18").Value = 0.05988023952095808
$ws.Range("S18").Value = 0.1137724550898204

# Row 19
$ws.Range("F19").Value = 0.01217765042979943
$ws.Range("H19").Value = 0.2041547277936963
$ws.Range("I19").Value = 0.0830945558739255
$ws.Range("J19").Value = 0.3911174785100286
$ws.Range("K19").Value = 0.1232091690544413
$ws.Range("M19").Value = 0.01790830945558739
$ws.Range("O19").Value = 0.08452722063037249
$ws.Range("S19").Value = 0.083810888252149
